# BurndownChart.20100927.xlsx — apply commit "Agregue los graficos de
# Costos, QC y Burndown Chart para entrega del 27/09"
#
# Net effect observed in the target OOXML:
#   - Sprint!A1:D1 header text is re-cased (title/weight/status/remaining
#     -> Title/Weight/Status/Remaining), which also renames the backing
#     table (Tabla2) columns.
#   - Sprint!A2:A12 task names are corrected (encoding fix for "n~" -> "ñ",
#     and a couple of rows re-worded) while the weight/status/remaining
#     values and formulas are untouched.
#   - 'Burndown Chart'!B2:B3 drop from 19 to 18 (the first two days of the
#     sprint burndown are corrected).
#   - The workbook now opens on the "Burndown Chart" tab instead of
#     "Sprint" (activeTab / tabSelected flip), and the last selection left
#     behind on the Sprint sheet is B2:B12.

$wb = $excel.ActiveWorkbook

$sprint = $wb.Worksheets.Item("Sprint")
$burndown = $wb.Worksheets.Item("Burndown Chart")

# --- Sprint: table header row -------------------------------------------
$sprint.Range("A1").Value = "Title"
$sprint.Range("B1").Value = "Weight"
$sprint.Range("C1").Value = "Status"
$sprint.Range("D1").Value = "Remaining"

# --- Sprint: task titles (column A), properly accented this time --------
$sprint.Range("A2").Value = "Crear VPC con ambiente de desarrollo"
$sprint.Range("A3").Value = "Crear el diseño general de la master page del sistema SelfManagement"
$sprint.Range("A4").Value = "Crear el mockup de la pagina de ABM de Campañas para los Jefes de Cuentas"
$sprint.Range("A5").Value = "Crear el mockup de la pagina de estadisticas globales de las Campañas para los Jefes de Cuentas (utilizando un dashboard y soportando busquedas)"
$sprint.Range("A6").Value = "Crear el mockup de la pagina de las estadisticas globales de un Supervisor para los Jefe de Cuentas y Supervisores (utilizando un dashboard)"
$sprint.Range("A7").Value = "Crear el mockup de la pagina de estadisticas y estado del sueldo variable de un Agente para los Jefes de Cuentas, Supervisores y Agentes (utilizar un dashboard)"
$sprint.Range("A8").Value = "Crear el mockup de la pagina de Login para los usuarios del sistema (Jefes de Cuentas, Supervisores y Agentes)"
$sprint.Range("A9").Value = "Crear el mockup de la pagina de ABM de usuario para el Responsable de IT"
$sprint.Range("A10").Value = "Crear la estructura inicial de la solucion SelfManagent con todos los proyectos requeridos"
$sprint.Range("A11").Value = "Diseñar el esquema de la base de datos para el sistema SelfManagement"
$sprint.Range("A12").Value = "Implementar la pantalla de alta de campañas para el sistema SelfManagement"

# --- Burndown Chart: corrected first two burndown points ----------------
$burndown.Range("B2").Value = 18
$burndown.Range("B3").Value = 18

# --- Leave the Sprint sheet selection where the author left it, then
#     switch focus to "Burndown Chart" (now the active tab on open) ------
$sprint.Activate()
$sprint.Range("B2:B12").Select()

$burndown.Activate()
